$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Quantity" column (column I),
# shifting Quantity/Price/Grant Date right by one.
$ws.Columns("I").Insert()
$ws.Columns("I").ColumnWidth = 14.8

# New header for the inserted column
$ws.Range("I1").Value = "Option Type"

# Populate the Option Type values for the two Options rows
$ws.Range("I6").Value = "Regular"
$ws.Range("I7").Value = "Phantom"

# Update the selection to match the post-edit active cell
$ws.Range("I7").Select()
